# Neural Network with very bad learning rate
# Rebuild the small "pipeline diagram" sheet with a third column (Fourier
# Spectrum) added to the left of the existing "High pass filtering" /
# "Machine Learning" blocks, each now a 2-column, 2-row sub-table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- start clean: drop old merges + contents so nothing from the old
#     D5:E5 / G5:H5 layout lingers once we move everything around ---
$ws.Cells.UnMerge()
$ws.Cells.Clear()

$xlCenter = -4108

# Row 4 - unchanged label
$ws.Range("A4").Value = "STA/LTA"

# Row 7 - the three section headers (each spans two columns, merged)
$ws.Range("A7").Value = "Fourier Spectrum"
$ws.Range("E7").Value = "High pass or band pass Filtering"
$ws.Range("H7").Value = "Machine Learning"

$ws.Range("A7:B7").HorizontalAlignment = $xlCenter
$ws.Range("A7:B7").WrapText = $true

$ws.Range("E7:F7").HorizontalAlignment = $xlCenter
$ws.Range("E7:F7").WrapText = $true

$ws.Range("H7:I7").HorizontalAlignment = $xlCenter
$ws.Range("H7:I7").VerticalAlignment = $xlCenter

$ws.Range("A7:B7").Merge()
$ws.Range("E7:F7").Merge()
$ws.Range("H7:I7").Merge()

# Row 8 - Input/Output sub-headers, bold, under every section
$ws.Range("A8").Value = "Input"
$ws.Range("B8").Value = "Output"
$ws.Range("E8").Value = "Input"
$ws.Range("F8").Value = "Output"
$ws.Range("H8").Value = "Input"
$ws.Range("I8").Value = "Output"
$ws.Range("A8:B8").Font.Bold = $true
$ws.Range("E8:F8").Font.Bold = $true
$ws.Range("H8:I8").Font.Bold = $true

# Row 9 - first data row for each section
$ws.Range("A9").Value = "Signal "
$ws.Range("B9").Value = "Frequencies"
$ws.Range("E9").Value = "Unfiltered Signal"
$ws.Range("F9").Value = "Filtered Signal"
$ws.Range("H9").Value = "Filtered signal"
$ws.Range("I9").Value = "Model"

# Row 10 - second data row (Fourier Spectrum + High pass filtering only)
$ws.Range("B10").Value = "Power/Amplitude"
$ws.Range("E10").Value = "Frequency intervals"
$ws.Range("H10").Value = "arrival time"

# Row 11 - Machine Learning output continues
$ws.Range("H11").Value = "End time"

# --- column widths (approximate the new layout; engine quantizes to the
#     nearest 1/6 character) ---
$ws.Columns.Item(1).ColumnWidth = 10.5
$ws.Columns.Item(2).ColumnWidth = 16.5
$ws.Columns.Item(4).ColumnWidth = 11.333333333333334
$ws.Columns.Item(5).ColumnWidth = 17.833333333333332
$ws.Columns.Item(6).ColumnWidth = 13

# --- view: zoom + selection ---
$ws.Application.ActiveWindow.Zoom = 175
$ws.Range("B14").Select() | Out-Null
